# ---------------------------------------------------------------------------
# Natmi following Dr Hou advice
#
# Refreshes the Adam10 -> Epha3 ligand-receptor edge table (Sheet1, A2:T13)
# with the revised NATMI statistics. The sending/target cluster combinations
# now cover ECs, FAPs and sCs as receiving clusters for every sending cluster
# (ECs, FAPs, M2, sCs), growing the table from 8 to 12 data rows.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order A..T:
#  Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
#  Ligand-expressing cells, Ligand detection rate,
#  Ligand average/total expression value,
#  Ligand derived specificity of average/total expression value,
#  Receptor-expressing cells, Receptor detection rate,
#  Receptor average/total expression value,
#  Receptor derived specificity of average/total expression value,
#  Edge average/total expression weight,
#  Edge average/total expression derived specificity
$rows = @(
    @("ECs","Adam10","Epha3","ECs",3.0,1.0,35.94400833333334,107.832025,0.3258225667258368,0.3258225667258368,2.0,0.6666666666666666,0.1220006666666667,0.366002,0.003000890193633104,0.003000890193633103,4.385192979338889,39.46673681405,0.000977757745351931,0.000977757745351931),
    @("ECs","Adam10","Epha3","FAPs",3.0,1.0,35.94400833333334,107.832025,0.3258225667258368,0.3258225667258368,3.0,1.0,40.12734033333333,120.382021,0.987025279393649,0.987025279393649,1442.337455335836,12981.03709802253,0.3215951099553249,0.3215951099553249),
    @("ECs","Adam10","Epha3","sCs",3.0,1.0,35.94400833333334,107.832025,0.3258225667258368,0.3258225667258368,3.0,1.0,0.4054843333333333,1.216453,0.009973830412717882,0.00997383041271788,14.57473225636945,131.172590307325,0.003249699025159952,0.003249699025159952),
    @("FAPs","Adam10","Epha3","ECs",3.0,1.0,23.289271,69.867813,0.2111108472940278,0.2111108472940278,2.0,0.6666666666666666,0.1220006666666667,0.366002,0.003000890193633104,0.003000890193633103,2.841306588180666,25.571759293626,0.0006335204714142236,0.0006335204714142236),
    @("FAPs","Adam10","Epha3","FAPs",3.0,1.0,23.289271,69.867813,0.2111108472940278,0.2111108472940278,3.0,1.0,40.12734033333333,120.382021,0.987025279393649,0.987025279393649,934.5365035322302,8410.828531790072,0.2083717430334177,0.2083717430334178),
    @("FAPs","Adam10","Epha3","sCs",3.0,1.0,23.289271,69.867813,0.2111108472940278,0.2111108472940278,3.0,1.0,0.4054843333333333,1.216453,0.009973830412717882,0.00997383041271788,9.443434525254332,84.99091072728899,0.002105583789195815,0.002105583789195814),
    @("M2","Adam10","Epha3","ECs",3.0,1.0,35.591128,106.773384,0.3226238033912781,0.3226238033912782,2.0,0.6666666666666666,0.1220006666666667,0.366002,0.003000890193633104,0.003000890193633103,4.342141343418666,39.07927209076799,0.000968158607829501,0.000968158607829501),
    @("M2","Adam10","Epha3","FAPs",3.0,1.0,35.591128,106.773384,0.3226238033912781,0.3226238033912782,3.0,1.0,40.12734033333333,120.382021,0.987025279393649,0.987025279393649,1428.177306103229,12853.59575492906,0.318437849681318,0.318437849681318),
    @("M2","Adam10","Epha3","sCs",3.0,1.0,35.591128,106.773384,0.3226238033912781,0.3226238033912782,3.0,1.0,0.4054843333333333,1.216453,0.009973830412717882,0.00997383041271788,14.43164480966133,129.884803286952,0.003217795102130644,0.003217795102130644),
    @("sCs","Adam10","Epha3","ECs",3.0,1.0,15.49333,46.47998999999999,0.1404427825888573,0.1404427825888573,2.0,0.6666666666666666,0.1220006666666667,0.366002,0.003000890193633104,0.003000890193633103,1.890196588886666,17.01176929998,0.0004214533690374479,0.0004214533690374479),
    @("sCs","Adam10","Epha3","FAPs",3.0,1.0,15.49333,46.47998999999999,0.1404427825888573,0.1404427825888573,3.0,1.0,40.12734033333333,120.382021,0.987025279393649,0.987025279393649,621.7061258066432,5595.355132259789,0.1386205767235884,0.1386205767235884),
    @("sCs","Adam10","Epha3","sCs",3.0,1.0,15.49333,46.47998999999999,0.1404427825888573,0.1404427825888573,3.0,1.0,0.4054843333333333,1.216453,0.009973830412717882,0.00997383041271788,6.282302586163333,56.54072327546999,0.00140075249623147,0.00140075249623147)
)

$numRows = $rows.Count
$numCols = 20
$data = New-Object "object[,]" $numRows,$numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $data[$i,$j] = $rows[$i][$j]
    }
}

$ws.Range("A2:T13").Value2 = $data
